$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 344, pushing the existing rows 344:355 down to 347:358
$ws.Rows("344:346").Insert()

# Populate the three newly inserted rows with the new weekly price entries
# Row 344 - Especial
$ws.Cells.Item(344, 1).Value = 8
$ws.Cells.Item(344, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(344, 3).Value = "Coquimbo"
$ws.Cells.Item(344, 4).Value = 45239
$ws.Cells.Item(344, 5).Value = 4
$ws.Cells.Item(344, 6).Value = "Fruta"
$ws.Cells.Item(344, 7).Value = 100107
$ws.Cells.Item(344, 8).Value = "Otros"
$ws.Cells.Item(344, 9).Value = 100107002
$ws.Cells.Item(344, 10).Value = "Chirimoya"
$ws.Cells.Item(344, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(344, 12).Value = "Especial"
$ws.Cells.Item(344, 13).Value = 300
$ws.Cells.Item(344, 14).Value = 15000
$ws.Cells.Item(344, 15).Value = 16000
$ws.Cells.Item(344, 16).Value = 15500
$ws.Cells.Item(344, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(344, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(344, 19).Value = 1550
$ws.Cells.Item(344, 20).Value = 10

# Row 345 - Primera
$ws.Cells.Item(345, 1).Value = 8
$ws.Cells.Item(345, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(345, 3).Value = "Coquimbo"
$ws.Cells.Item(345, 4).Value = 45239
$ws.Cells.Item(345, 5).Value = 4
$ws.Cells.Item(345, 6).Value = "Fruta"
$ws.Cells.Item(345, 7).Value = 100107
$ws.Cells.Item(345, 8).Value = "Otros"
$ws.Cells.Item(345, 9).Value = 100107002
$ws.Cells.Item(345, 10).Value = "Chirimoya"
$ws.Cells.Item(345, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(345, 12).Value = "Primera"
$ws.Cells.Item(345, 13).Value = 360
$ws.Cells.Item(345, 14).Value = 11000
$ws.Cells.Item(345, 15).Value = 12000
$ws.Cells.Item(345, 16).Value = 11500
$ws.Cells.Item(345, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(345, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(345, 19).Value = 1150
$ws.Cells.Item(345, 20).Value = 10

# Row 346 - Segunda
$ws.Cells.Item(346, 1).Value = 8
$ws.Cells.Item(346, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(346, 3).Value = "Coquimbo"
$ws.Cells.Item(346, 4).Value = 45239
$ws.Cells.Item(346, 5).Value = 4
$ws.Cells.Item(346, 6).Value = "Fruta"
$ws.Cells.Item(346, 7).Value = 100107
$ws.Cells.Item(346, 8).Value = "Otros"
$ws.Cells.Item(346, 9).Value = 100107002
$ws.Cells.Item(346, 10).Value = "Chirimoya"
$ws.Cells.Item(346, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(346, 12).Value = "Segunda"
$ws.Cells.Item(346, 13).Value = 240
$ws.Cells.Item(346, 14).Value = 8000
$ws.Cells.Item(346, 15).Value = 9000
$ws.Cells.Item(346, 16).Value = 8500
$ws.Cells.Item(346, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(346, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(346, 19).Value = 850
$ws.Cells.Item(346, 20).Value = 10
